# Update the "想去人数" (F column) counts on the 展览, 演出 and 全部类型 sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 5372
$ws1.Range("F9").Value  = 280
$ws1.Range("F12").Value = 2647
$ws1.Range("F13").Value = 2647
$ws1.Range("F15").Value = 2354
$ws1.Range("F16").Value = 2354
$ws1.Range("F26").Value = 168
$ws1.Range("F27").Value = 106
$ws1.Range("F43").Value = 73
$ws1.Range("F47").Value = 573
$ws1.Range("F49").Value = 113

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 23
$ws2.Range("F5").Value = 224

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5372
$ws4.Range("F7").Value  = 280
$ws4.Range("F10").Value = 2647
$ws4.Range("F12").Value = 23
$ws4.Range("F13").Value = 224
$ws4.Range("F15").Value = 2354
$ws4.Range("F28").Value = 168
$ws4.Range("F29").Value = 106
$ws4.Range("F42").Value = 73
$ws4.Range("F47").Value = 113
